{"js": "// Remove non-breaking spaces (U+00A0) from the digest template.\n// Rule (per commit message): either remove the NBSP if it sits next to a\n// regular space already, or otherwise substitute it with a regular space.\nconst body = context.document.body;\n\n// 1) \"...errors in the<NBSP>tests\" -> \"...errors in the tests\"\n//    (no adjacent space, so the NBSP is substituted with a regular space)\nconst r1 = body.search(\"the\\u00a0tests\", { matchCase: true, matchWildcards: false });\ncontext.load(r1, \"text\");\n\n// 2) \"...we've used <NBSP>before...\" -> \"...we've used before...\"\n//    (already preceded by a regular space, so the NBSP is simply removed)\nconst r2 = body.search(\" \\u00a0before\", { matchCase: true, matchWildcards: false });\ncontext.load(r2, \"text\");\n\n// 3) The bold run \"It's not just mammals who can recognise sample data.\"\n//    gains a trailing regular space (the leading NBSP that used to open the\n//    following run is normalised away, and its \"space\" ends up here).\nconst r3 = body.search(\"who can recognise sample data.\", { matchCase: true, matchWildcards: false });\ncontext.load(r3, \"text\");\n\n// 4) \"<NBSP>Image credit:<NBSP>Anonymous and Anonymous<NBSP>(CC BY<NBSP>4.0)\"\n//    -> \"Image credit: Anonymous and Anonymous (CC BY 4.0)\"\n//    (leading NBSP removed entirely - see point 3 - remaining NBSPs become spaces)\nconst r4 = body.search(\n  \"\\u00a0Image credit:\\u00a0Anonymous and Anonymous\\u00a0(CC BY\\u00a04.0)\",\n  { matchCase: true, matchWildcards: false }\n);\ncontext.load(r4, \"text\");\n\nawait context.sync();\n\nif (r1.items.length === 1) {\n  r1.items[0].insertText(\"the tests\", Word.InsertLocation.replace);\n}\nif (r2.items.length === 1) {\n  r2.items[0].insertText(\" before\", Word.InsertLocation.replace);\n}\nif (r3.items.length === 1) {\n  // Append - keeps this text in its own (bold) run rather than merging\n  // with the differently-formatted run that follows it.\n  r3.items[0].insertText(\" \", Word.InsertLocation.end);\n}\nif (r4.items.length === 1) {\n  r4.items[0].insertText(\n    \"Image credit: Anonymous and Anonymous (CC BY 4.0)\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Remove non-breaking spaces (U+00A0) from the digest template.\n# Rule (per commit message): either remove the NBSP if it sits next to a\n# regular space already, or otherwise substitute it with a regular space.\n\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # wdFindWrap=1 (wdFindContinue) via the Wrap slot is passed positionally below;\n    # last arg 1 = wdReplaceOne so only the located occurrence is touched.\n    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    return $found\n}\n\n# 1) \"...errors in the<NBSP>tests\" -> \"...errors in the tests\"\n#    (no adjacent space, so the NBSP is substituted with a regular space)\nReplace-Once \"the${nbsp}tests\" \"the tests\" | Out-Null\n\n# 2) \"...we've used <NBSP>before...\" -> \"...we've used before...\"\n#    (already preceded by a regular space, so the NBSP is simply removed)\nReplace-Once \" ${nbsp}before\" \" before\" | Out-Null\n\n# 3) The bold run \"It's not just mammals who can recognise sample data.\"\n#    gains a trailing regular space (the leading NBSP that used to open the\n#    following run is normalised away - see step 4 - and its \"space\" ends\n#    up here instead).\nReplace-Once \"who can recognise sample data.\" \"who can recognise sample data. \" | Out-Null\n\n# 4) \"<NBSP>Image credit:<NBSP>Anonymous and Anonymous<NBSP>(CC BY<NBSP>4.0)\"\n#    -> \"Image credit: Anonymous and Anonymous (CC BY 4.0)\"\n#    (leading NBSP removed entirely - see step 3 - remaining NBSPs become spaces)\nReplace-Once \"${nbsp}Image credit:${nbsp}Anonymous and Anonymous${nbsp}(CC BY${nbsp}4.0)\" \"Image credit: Anonymous and Anonymous (CC BY 4.0)\" | Out-Null\n"}
